$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new quote row right after the current last row (row 79 -> row 80).
# The leading apostrophe forces the date-looking string in column A to be
# stored as plain text instead of being auto-parsed into a date serial value,
# matching the existing text values used throughout the sheet.
$ws.Range("A80").Value = "'2025-10-15"
$ws.Range("B80").Value = "21:22:21"
$ws.Range("C80").Value = "1.00 EUR = 1,675.9082"
